$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3200.3333
$ws.Range("I2").Value = 301
$ws.Range("J2").Value = 4650
$ws.Range("K2").Value = 301
$ws.Range("L2").Value = 4650
$ws.Range("M2").Value = -188
$ws.Range("N2").Value = -4876
$ws.Range("H12").Value = 919.53845
$ws.Range("I12").Value = 914.9524
$ws.Range("J12").Value = 938.8
$ws.Range("K12").Value = 914.9524
$ws.Range("L12").Value = 938.8
$ws.Range("M12").Value = -744.9524
$ws.Range("N12").Value = -1278.8
$ws.Range("H18").Value = 197.5
$ws.Range("I18").Value = 197.5
$ws.Range("K18").Value = 197.5
$ws.Range("M18").Value = 86.5
$ws.Range("H19").Value = 1940.0526
$ws.Range("I19").Value = 1049.875
$ws.Range("J19").Value = 2587.4546
$ws.Range("K19").Value = 1049.875
$ws.Range("L19").Value = 2587.4546
$ws.Range("M19").Value = -874.875
$ws.Range("N19").Value = -2937.4546
$ws.Range("H40").Value = 8530.4
$ws.Range("J40").Value = 8530.4
$ws.Range("L40").Value = 8530.4
$ws.Range("N40").Value = -8880.4
$ws.Range("H55").Value = 164.76923
$ws.Range("J55").Value = 263.85715
$ws.Range("L55").Value = 263.85715
$ws.Range("N55").Value = -691.85715
$ws.Range("H74").Value = 17863186
$ws.Range("I74").Value = 20413284
$ws.Range("K74").Value = 20413284
$ws.Range("M74").Value = -20412348
$ws.Range("H77").Value = 17863186
$ws.Range("I77").Value = 20413284
$ws.Range("K77").Value = 102066420
$ws.Range("M77").Value = -102061740
$ws.Range("H112").Value = 2373.425
$ws.Range("J112").Value = 2486.1353
$ws.Range("L112").Value = 7458.4059
$ws.Range("N112").Value = -9674.4059
$ws.Range("H129").Value = 1247.75
$ws.Range("I129").Value = 997
$ws.Range("K129").Value = 2991
$ws.Range("M129").Value = 2009
$ws.Range("H132").Value = 507010.66
$ws.Range("I132").Value = 567720.4
$ws.Range("K132").Value = 1703161.2
$ws.Range("M132").Value = -1700631.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1879.8462
$ws.Range("J63").Value = 1999.5
$ws.Range("L63").Value = 1999.5
$ws.Range("N63").Value = -3371.5
$ws.Range("H66").Value = 1879.8462
$ws.Range("J66").Value = 1999.5
$ws.Range("L66").Value = 9997.5
$ws.Range("N66").Value = -16861.5
$ws.Range("H110").Value = 5667.6665
$ws.Range("I110").Value = 2772.125
$ws.Range("J110").Value = 7115.4375
$ws.Range("K110").Value = 2772.125
$ws.Range("L110").Value = 7115.4375
$ws.Range("M110").Value = -727.125
$ws.Range("N110").Value = -11205.4375
$ws.Range("H122").Value = 2723.423
$ws.Range("I122").Value = 1827.8096
$ws.Range("J122").Value = 6485
$ws.Range("K122").Value = 5483.4288
$ws.Range("L122").Value = 19455
$ws.Range("M122").Value = -3033.4288
$ws.Range("N122").Value = -24355

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 992.25
$ws.Range("I22").Value = 485
$ws.Range("J22").Value = 1499.5
$ws.Range("K22").Value = 485
$ws.Range("L22").Value = 1499.5
$ws.Range("M22").Value = -312
$ws.Range("N22").Value = -1845.5
$ws.Range("H134").Value = 1028956.2
$ws.Range("I134").Value = 1094752.2
$ws.Range("J134").Value = 20083
$ws.Range("K134").Value = 3284256.6
$ws.Range("L134").Value = 60249
$ws.Range("M134").Value = -3281721.6
$ws.Range("N134").Value = -65319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 15155214
$ws.Range("I16").Value = 33335180
$ws.Range("K16").Value = 33335180
$ws.Range("M16").Value = -33334893
$ws.Range("H31").Value = 6659.5264
$ws.Range("I31").Value = 1105.7778
$ws.Range("J31").Value = 11657.9
$ws.Range("K31").Value = 1105.7778
$ws.Range("L31").Value = 11657.9
$ws.Range("M31").Value = -810.7778000000001
$ws.Range("N31").Value = -12247.9
$ws.Range("H32").Value = 5875.25
$ws.Range("I32").Value = 755
$ws.Range("J32").Value = 10995.5
$ws.Range("K32").Value = 755
$ws.Range("L32").Value = 10995.5
$ws.Range("M32").Value = -439
$ws.Range("N32").Value = -11627.5
$ws.Range("H34").Value = 6659.5264
$ws.Range("I34").Value = 1105.7778
$ws.Range("J34").Value = 11657.9
$ws.Range("K34").Value = 1105.7778
$ws.Range("L34").Value = 11657.9
$ws.Range("M34").Value = -903.7778000000001
$ws.Range("N34").Value = -12061.9
$ws.Range("H99").Value = 13892227
$ws.Range("I99").Value = 15875971
$ws.Range("K99").Value = 15875971
$ws.Range("M99").Value = -15874473
$ws.Range("H113").Value = 15155214
$ws.Range("I113").Value = 33335180
$ws.Range("K113").Value = 33335180
$ws.Range("M113").Value = -33333010
$ws.Range("H126").Value = 13892227
$ws.Range("I126").Value = 15875971
$ws.Range("K126").Value = 47627913
$ws.Range("M126").Value = -47625443
$ws.Range("H132").Value = 7001.8477
$ws.Range("I132").Value = 4493.973
$ws.Range("J132").Value = 17312
$ws.Range("K132").Value = 13481.919
$ws.Range("L132").Value = 51936
$ws.Range("M132").Value = -10951.919
$ws.Range("N132").Value = -56996
$ws.Range("H134").Value = 90917336
$ws.Range("I134").Value = 111117330
$ws.Range("J134").Value = 17375
$ws.Range("K134").Value = 333351990
$ws.Range("L134").Value = 52125
$ws.Range("M134").Value = -333349455
$ws.Range("N134").Value = -57195
$ws.Range("H141").Value = 432499.34
$ws.Range("J141").Value = 487999.8
$ws.Range("L141").Value = 487999.8
$ws.Range("N141").Value = -498359.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 244.0625
$ws.Range("I26").Value = 239.58333
$ws.Range("J26").Value = 257.5
$ws.Range("K26").Value = 718.74999
$ws.Range("L26").Value = 772.5
$ws.Range("M26").Value = -430.74999
$ws.Range("N26").Value = -1348.5
$ws.Range("H82").Value = 16306
$ws.Range("I82").Value = 6688.75
$ws.Range("J82").Value = 23999.8
$ws.Range("K82").Value = 20066.25
$ws.Range("L82").Value = 71999.39999999999
$ws.Range("M82").Value = -19660.25
$ws.Range("N82").Value = -72811.39999999999
$ws.Range("H85").Value = 16306
$ws.Range("I85").Value = 6688.75
$ws.Range("J85").Value = 23999.8
$ws.Range("K85").Value = 20066.25
$ws.Range("L85").Value = 71999.39999999999
$ws.Range("M85").Value = -18662.25
$ws.Range("N85").Value = -74807.39999999999
$ws.Range("H98").Value = 1358156.6
$ws.Range("J98").Value = 1832.5
$ws.Range("L98").Value = 5497.5
$ws.Range("N98").Value = -8493.5
$ws.Range("H134").Value = 9389.956
$ws.Range("I134").Value = 8453.137000000001
$ws.Range("K134").Value = 25359.411
$ws.Range("M134").Value = -20289.411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 870331.1
$ws.Range("I102").Value = 1555483.4
$ws.Range("J102").Value = 6443.4346
$ws.Range("K102").Value = 1555483.4
$ws.Range("L102").Value = 6443.4346
$ws.Range("M102").Value = -1553861.4
$ws.Range("N102").Value = -9687.434600000001
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H122").Value = 5018.0967
$ws.Range("I122").Value = 4021.8
$ws.Range("K122").Value = 12065.4
$ws.Range("M122").Value = -9615.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8189.45
$ws.Range("I7").Value = 8254.916999999999
$ws.Range("K7").Value = 8254.916999999999
$ws.Range("M7").Value = -8142.916999999999
$ws.Range("H46").Value = 45455610
$ws.Range("I46").Value = 1099.5714
$ws.Range("J46").Value = 125001000
$ws.Range("K46").Value = 1099.5714
$ws.Range("L46").Value = 125001000
$ws.Range("M46").Value = -911.5714
$ws.Range("N46").Value = -125001376
$ws.Range("H82").Value = 5756.9287
$ws.Range("I82").Value = 2119.9
$ws.Range("K82").Value = 2119.9
$ws.Range("M82").Value = -1758.9
$ws.Range("H85").Value = 5756.9287
$ws.Range("I85").Value = 2119.9
$ws.Range("K85").Value = 2119.9
$ws.Range("M85").Value = -871.9000000000001
$ws.Range("H122").Value = 3474.1428
$ws.Range("I122").Value = 2733.9048
$ws.Range("K122").Value = 8201.714399999999
$ws.Range("M122").Value = -5751.714399999999
$ws.Range("H126").Value = 8189.45
$ws.Range("I126").Value = 8254.916999999999
$ws.Range("K126").Value = 24764.751
$ws.Range("M126").Value = -22294.751
$ws.Range("H132").Value = 9617.344999999999
$ws.Range("I132").Value = 8276.727999999999
$ws.Range("K132").Value = 24830.184
$ws.Range("M132").Value = -22300.184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8334316.5
$ws.Range("I107").Value = 12500710
$ws.Range("J107").Value = 1530
$ws.Range("K107").Value = 37502130
$ws.Range("L107").Value = 4590
$ws.Range("M107").Value = -37500210
$ws.Range("N107").Value = -8430
